$wb = $excel.ActiveWorkbook

# --- NewLoanInput sheet: rename loan product, insert "Firstrepaymenton" row ---
$ws1 = $wb.Worksheets.Item("NewLoanInput")

# "Chaithanya 123" -> "chaithanyatest"
$ws1.Range("B2").Value = "chaithanyatest"

# Insert a new row 7 for "Firstrepaymenton" (shifts old rows 7-20 down to 8-21)
$ws1.Rows("7").Insert()
$ws1.Range("A7").Value = "Firstrepaymenton"
$ws1.Range("B7").Value = "2/1/2015"

$ws1.Range("B7").Select()

# --- Summary sheet: selection only changes ---
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("A4").Select()

# --- Repayment Schedule sheet: selection only changes ---
$ws3 = $wb.Worksheets.Item("Repayment Schedule")
$ws3.Range("F6").Select()

# --- Transactions sheet: selection + Entry ID value change ---
$ws4 = $wb.Worksheets.Item("Transactions")
$ws4.Range("A2").Value = 227
$ws4.Range("G2").Select()

Write-Output "done"
